$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Printed On" date bump ---
$ws.Range("M3").Value = "Printed On: 10/23/2025"

# --- Agency name change (keep "Agency:" bold, city name non-bold, like the source) ---
$agencyCell = $ws.Range("B16")
$agencyCell.Characters(8, 14).Text = " MCKEES ROCKS BOROUGH"
$agencyRun1 = $agencyCell.Characters(1, 7)
$agencyRun1.Font.Bold = $true
$agencyRun1.Font.Italic = $false
$agencyRun1.Font.Underline = $false
$agencyRun1.Font.Strikethrough = $false
$agencyRun2 = $agencyCell.Characters(8, 22)
$agencyRun2.Font.Bold = $false
$agencyRun2.Font.Italic = $false
$agencyRun2.Font.Underline = $false
$agencyRun2.Font.Strikethrough = $false

# --- Updated Return A offense counts ---
$ws.Range("E21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("E24").Value = 9
$ws.Range("I24").Value = 9
$ws.Range("J24").Value = 3
$ws.Range("N24").Value = 1
$ws.Range("E25").Value = 5
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 1
$ws.Range("N25").Value = 1
$ws.Range("E28").Value = 4
$ws.Range("I28").Value = 4
$ws.Range("E29").Value = 17
$ws.Range("I29").Value = 17
$ws.Range("J29").Value = 9
$ws.Range("N29").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("E31").Value = 2
$ws.Range("I31").Value = 2
$ws.Range("J31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("E32").Value = 3
$ws.Range("I32").Value = 3
$ws.Range("J32").Value = 1
$ws.Range("E33").Value = 9
$ws.Range("I33").Value = 9
$ws.Range("J33").Value = 8
$ws.Range("E34").Value = 24
$ws.Range("I34").Value = 24
$ws.Range("J34").Value = 3
$ws.Range("N34").Value = 0
$ws.Range("E35").Value = 14
$ws.Range("I35").Value = 14
$ws.Range("J35").Value = 1
$ws.Range("N35").Value = 0
$ws.Range("E36").Value = 9
$ws.Range("I36").Value = 9
$ws.Range("J36").Value = 2
$ws.Range("E38").Value = 171
$ws.Range("I38").Value = 171
$ws.Range("J38").Value = 33
$ws.Range("N38").Value = 4
$ws.Range("E39").Value = 34
$ws.Range("I39").Value = 34
$ws.Range("J39").Value = 5
$ws.Range("N39").Value = 1
$ws.Range("E40").Value = 33
$ws.Range("I40").Value = 33
$ws.Range("J40").Value = 5
$ws.Range("N40").Value = 1
$ws.Range("E41").Value = 1
$ws.Range("I41").Value = 1
$ws.Range("E44").Value = 1
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 0
$ws.Range("E45").Value = 256
$ws.Range("I45").Value = 256
$ws.Range("J45").Value = 53
$ws.Range("N45").Value = 6
$ws.Range("E46").Value = 132
$ws.Range("I46").Value = 132
$ws.Range("J46").Value = 80
$ws.Range("N46").Value = 9
$ws.Range("E48").Value = 16
$ws.Range("I48").Value = 16
$ws.Range("E50").Value = 7
$ws.Range("I50").Value = 7
$ws.Range("J50").Value = 3
$ws.Range("N50").Value = 2
$ws.Range("E51").Value = 88
$ws.Range("I51").Value = 88
$ws.Range("J51").Value = 8
$ws.Range("E52").Value = 22
$ws.Range("I52").Value = 22
$ws.Range("J52").Value = 14
$ws.Range("N52").Value = 1
$ws.Range("E53").Value = 1
$ws.Range("I53").Value = 1
$ws.Range("E54").Value = 3
$ws.Range("I54").Value = 3
$ws.Range("E55").Value = 43
$ws.Range("I55").Value = 43
$ws.Range("J55").Value = 29
$ws.Range("J56").Value = 2
$ws.Range("J57").Value = 2
$ws.Range("E61").Value = 40
$ws.Range("I61").Value = 40
$ws.Range("J61").Value = 27
$ws.Range("E62").Value = 20
$ws.Range("I62").Value = 20
$ws.Range("J62").Value = 11
$ws.Range("E63").Value = 15
$ws.Range("I63").Value = 15
$ws.Range("J63").Value = 11
$ws.Range("E65").Value = 4
$ws.Range("I65").Value = 4
$ws.Range("J65").Value = 4
$ws.Range("E71").Value = 12
$ws.Range("I71").Value = 12
$ws.Range("J71").Value = 12
$ws.Range("E73").Value = 2
$ws.Range("I73").Value = 2
$ws.Range("J73").Value = 2
$ws.Range("E74").Value = 10
$ws.Range("I74").Value = 10
$ws.Range("J74").Value = 10
$ws.Range("E75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("E76").Value = 45
$ws.Range("I76").Value = 45
$ws.Range("J76").Value = 44
$ws.Range("N76").Value = 4
$ws.Range("E77").Value = 383
$ws.Range("I77").Value = 383
$ws.Range("J77").Value = 206
$ws.Range("N77").Value = 16
